$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose pk0 values are
# "2eb86e3a-3239-440b-8e9f-7ab13604494d" and "5434643f-0598-4c31-833d-d759fdc14e84"
# (originally rows 4 and 5). This shifts subsequent rows up.
$ws.Rows("4:5").Delete()

# Update the D (dt1) and E (dt2) columns for the remaining data rows to their
# new values (E values rounded to 9 decimal places).
$ws.Cells.Item(2, 4).Value2 = 4
$ws.Cells.Item(2, 5).Value2 = 7.412738474

$ws.Cells.Item(3, 4).Value2 = 10
$ws.Cells.Item(3, 5).Value2 = 3.109741954

$ws.Cells.Item(4, 4).Value2 = 9
$ws.Cells.Item(4, 5).Value2 = 10.49718063

$ws.Cells.Item(5, 4).Value2 = 6
$ws.Cells.Item(5, 5).Value2 = 7.906586852

$ws.Cells.Item(6, 4).Value2 = 3
$ws.Cells.Item(6, 5).Value2 = 8.288637606

$ws.Cells.Item(7, 4).Value2 = 7
$ws.Cells.Item(7, 5).Value2 = 10.45706548

$ws.Cells.Item(8, 4).Value2 = 8
$ws.Cells.Item(8, 5).Value2 = 8.264160006

$ws.Cells.Item(9, 4).Value2 = 6
$ws.Cells.Item(9, 5).Value2 = 8.182236858

# Update the selected cell to match the saved view state.
$ws.Range("K7").Select()
